$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns AD:AF hold the team's season record (Wins / Losses / Ties),
# mirroring the style already used by the other header cells (e.g. AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row gets the same team season record: 87-75-0.
for ($r = 2; $r -le 58; $r++) {
    $ws.Cells.Item($r, 30).Value = 87
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
